$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45174
}
